$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.411.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7172"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.72%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07971"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3151"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("E11").Value = "  -3.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.893.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.242"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7079"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.419"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008439"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.410.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.136.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.681"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1589"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.073"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("E28").Value = "  +1.84%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.422"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.312"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.223"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05330"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.952"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7570"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.177"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01896"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.274.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.66%  "

$ws.Range("E40").Value = "  +0.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.76%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9062"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.57%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000129"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.030.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.808"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5203"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.528"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4348"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
